$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new values look like plain numbers as Text,
# so Excel stores them as strings (matching the source data which uses
# formatted/rounded price strings, not numeric values).
$textCellsRange = $ws.Range("D5,D6,D7,D8,D9,D10,D11,D14,D15,D17,D19,D20,D21,D22,D24,D25,D26,D29,D30,D31,D32,D33,D34,D35,D36,D38,D39,D41,D42,D44,D46,D47,D48,D49,D50,D51")
foreach ($area in $textCellsRange.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range("D2").Value = "57.391.11"
$ws.Range("E2").Value = "  +1.29%  "
$ws.Range("D3").Value = "3.011.67"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "507.22"
$ws.Range("D6").Value = "139.75"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "0.434"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "7.57"
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").Value = "0.110"
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("D11").Value = "0.366"
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("D12").Value = "3.530.53"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").Value = "26.29"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("D15").Value = "0.0000162"
$ws.Range("E15").Value = "  +3.45%  "
$ws.Range("D16").Value = "57.446.07"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "6.19"
$ws.Range("E17").Value = "  +4.04%  "
$ws.Range("D18").Value = "3.011.72"
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("D19").Value = "12.85"
$ws.Range("E19").Value = "  +2.75%  "
$ws.Range("D20").Value = "7.95"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "327.24"
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "0.995"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").Value = "0.498"
$ws.Range("E24").Value = "  +2.79%  "
$ws.Range("D25").Value = "64.58"
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("D26").Value = "0.166"
$ws.Range("E26").Value = "  -3.09%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "0.0₃0921"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("D29").Value = "6.74"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("D30").Value = "7.34"
$ws.Range("E30").Value = "  +5.37%  "
$ws.Range("D31").Value = "1.80"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").Value = "20.59"
$ws.Range("E33").Value = "  -0.30%  "
$ws.Range("D34").Value = "4.76"
$ws.Range("E34").Value = "  +5.00%  "
$ws.Range("D35").Value = "153.61"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").Value = "5.89"
$ws.Range("E37").Value = "  +0.68%  "
$ws.Range("D38").Value = "24.60"
$ws.Range("E38").Value = "  +3.84%  "
$ws.Range("D39").Value = "0.0677"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").Value = "3.043.45"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").Value = "37.82"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").Value = "3.85"
$ws.Range("E42").Value = "  +4.53%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "0.649"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.229.08"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "1.41"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "0.975"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("D48").Value = "6.07"
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("D49").Value = "0.0239"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "19.57"
$ws.Range("E50").Value = "  +1.55%  "
$ws.Range("D51").Value = "1.86"
$ws.Range("E51").Value = "  -4.29%  "
